$wb = $excel.ActiveWorkbook

# --- "Hands" sheet: add Hand9 and Hand10 ---
$hands = $wb.Worksheets.Item("Hands")
$hands.Range("A9").Value = "Hand9"
$hands.Range("B9").Value = "Eoghan Mac Gilleoin"
$hands.Range("A10").Value = "Hand10"
$hands.Range("B10").Value = "Niall Mac Mhuirich"
$hands.Range("C9").Value = "Transcription 9"
$hands.Range("C10").Value = "Transcription 12"
[void]$hands.Range("C17").Select()

# --- "Glyphs" sheet: add g20 "Superscript e", rename "ir superscript" -> "Superscript ir" ---
$glyphs = $wb.Worksheets.Item("Glyphs")
$glyphs.Range("A21").Value = "g20"
$glyphs.Range("B21").Value = "Superscript e"
$glyphs.Range("B7").Value = "Superscript ir"
[void]$glyphs.Range("B7").Select()
